# Hortaliza, Macroferia Regional de Talca - Repollo
# Insert two new weekly price rows (date serial 45194) right before the
# existing row 508, shifting all subsequent rows down by two (508-556 -> 510-558).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 508-509; everything from the old row 508 onward
# (previously ending at 556) moves down to 510-558, and the sheet's used
# range grows accordingly (dimension becomes A1:R558).
$ws.Rows("508:509").Insert()

# New row 508: Primera quality entry for date 45194
$ws.Range("A508").Value2 = 5
$ws.Range("B508").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C508").Value2 = 'Maule'
$ws.Range("D508").Value2 = 45194
$ws.Range("E508").Value2 = 7
$ws.Range("F508").Value2 = 100112006
$ws.Range("G508").Value2 = 'Repollo'
$ws.Range("H508").Value2 = 'Crespo record'
$ws.Range("I508").Value2 = 'Primera'
$ws.Range("J508").Value2 = 2000
$ws.Range("K508").Value2 = 900
$ws.Range("L508").Value2 = 900
$ws.Range("M508").Value2 = 900
$ws.Range("N508").Value2 = '$/unidad'
$ws.Range("O508").Value2 = 'Región del Maule'
$ws.Range("P508").Value2 = 900
$ws.Range("Q508").Value2 = 1
$ws.Range("R508").Value2 = 'Hortaliza'

# New row 509: Segunda quality entry for the same date 45194
$ws.Range("A509").Value2 = 5
$ws.Range("B509").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C509").Value2 = 'Maule'
$ws.Range("D509").Value2 = 45194
$ws.Range("E509").Value2 = 7
$ws.Range("F509").Value2 = 100112006
$ws.Range("G509").Value2 = 'Repollo'
$ws.Range("H509").Value2 = 'Crespo record'
$ws.Range("I509").Value2 = 'Segunda'
$ws.Range("J509").Value2 = 3000
$ws.Range("K509").Value2 = 500
$ws.Range("L509").Value2 = 500
$ws.Range("M509").Value2 = 500
$ws.Range("N509").Value2 = '$/unidad'
$ws.Range("O509").Value2 = 'Región del Maule'
$ws.Range("P509").Value2 = 500
$ws.Range("Q509").Value2 = 1
$ws.Range("R509").Value2 = 'Hortaliza'
